$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LEAVE CREDITS")

# Header area: F4 gets a new entry "TICC"
$ws.Range("F4").Value = "TICC"

# Monthly leave-credit rows 12-20: add the period date (month-end) in column A
$ws.Range("A12").Value = 45077
$ws.Range("A13").Value = 45107
$ws.Range("A14").Value = 45138
$ws.Range("A15").Value = 45169
$ws.Range("A16").Value = 45199
$ws.Range("A17").Value = 45230
$ws.Range("A18").Value = 45260
$ws.Range("A19").Value = 45291
$ws.Range("A20").Value = 45322

# EARNED column (C) for rows 12-18 : 1.25 credits per month
$ws.Range("C12").Value = 1.25
$ws.Range("C13").Value = 1.25
$ws.Range("C14").Value = 1.25
$ws.Range("C15").Value = 1.25
$ws.Range("C16").Value = 1.25
$ws.Range("C17").Value = 1.25
$ws.Range("C18").Value = 1.25

# Remarks / particulars text on rows 17-18
$ws.Range("B17").Value = "SP(1-0-0)"
$ws.Range("B18").Value = "SL(2-0-0)"

# Absence undertime w/ pay value for row 18
$ws.Range("H18").Value = 2

# Extra remark date/text in column K
$ws.Range("K17").Value = 45205
$ws.Range("K17").NumberFormat = "mm-dd-yy"
$ws.Range("K18").Value = "11/5,19/2023"

$wb.Save()
